$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4327
$ws1.Range("F3").Value = 2442
$ws1.Range("F6").Value = 45
$ws1.Range("F10").Value = 145
$ws1.Range("F12").Value = 1606
$ws1.Range("F13").Value = 294
$ws1.Range("F14").Value = 3389

# Sheet "演出" (sheet2): update "想去人数" (F column) values
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 3

# Sheet "全部类型" (sheet4): update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4327
$ws4.Range("F3").Value = 2442
$ws4.Range("F7").Value = 45
$ws4.Range("F12").Value = 145
$ws4.Range("F14").Value = 3
$ws4.Range("F16").Value = 1606
$ws4.Range("F17").Value = 294
$ws4.Range("F18").Value = 3389

$wb.Save()
